$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.05107311783547175
$ws.Range("D2").Value = 0.1348425512076261
$ws.Range("E2").Value = 0.07753440865865002
$ws.Range("F2").Value = 7.336283575130153
$ws.Range("G2").Value = 0.002807598003919605
$ws.Range("I2").Value = 5.48723525839327
$ws.Range("J2").Value = 0.1669609229091691
$ws.Range("L2").Value = 0.1653712041401434
$ws.Range("M2").Value = 12.85679754728847

$ws.Range("C3").Value = 0.04567767352035901
$ws.Range("D3").Value = 0.1290583356567936
$ws.Range("E3").Value = 0.06836154237538494
$ws.Range("F3").Value = 7.535399107722242
$ws.Range("G3").Value = 0.002822979187190019
$ws.Range("I3").Value = 5.635235392315465
$ws.Range("J3").Value = 0.1617416232217437
$ws.Range("L3").Value = 0.1497844912632473
$ws.Range("M3").Value = 11.66965917421413

$ws.Range("C4").Value = 0.04240839678877251
$ws.Range("D4").Value = 0.1257111067537693
$ws.Range("E4").Value = 0.06273122383467467
$ws.Range("F4").Value = 7.66756436280896
$ws.Range("G4").Value = 0.002832821965852642
$ws.Range("I4").Value = 5.733589757468593
$ws.Range("J4").Value = 0.1587253388242544
$ws.Range("L4").Value = 0.1402753348585151
$ws.Range("M4").Value = 10.94011934200449

$ws.Range("C5").Value = 0.04108640659902107
$ws.Range("D5").Value = 0.1243967398262527
$ws.Range("E5").Value = 0.06043684315277886
$ws.Range("F5").Value = 7.72388160468725
$ws.Range("G5").Value = 0.002836934067618377
$ws.Range("I5").Value = 5.775526339184779
$ws.Range("J5").Value = 0.1575422248114222
$ws.Range("L5").Value = 0.1364145340559304
$ws.Range("M5").Value = 10.64260377035265

$ws.Range("C6").Value = 0.04086749169991322
$ws.Range("D6").Value = 0.1241814401303998
$ws.Range("E6").Value = 0.06005585180788842
$ws.Range("F6").Value = 7.733380613111891
$ws.Range("G6").Value = 0.002837623008571925
$ws.Range("I6").Value = 5.782601280452852
$ws.Range("J6").Value = 0.1573485136467383
$ws.Range("L6").Value = 0.1357742811729423
$ws.Range("M6").Value = 10.59318636903942

$ws.Range("C7").Value = 0.04239052728328829
$ws.Range("D7").Value = 0.1256931819150537
$ws.Range("E7").Value = 0.06270028149939577
$ws.Range("F7").Value = 7.668313963369201
$ws.Range("G7").Value = 0.002832877012571423
$ws.Range("I7").Value = 5.73414784489681
$ws.Range("J7").Value = 0.1587091981814268
$ws.Range("L7").Value = 0.1402232104336747
$ws.Range("M7").Value = 10.93610792236916

$ws.Range("C8").Value = 0.04920330937341078
$ws.Range("D8").Value = 0.1328046978263728
$ws.Range("E8").Value = 0.07437093538643325
$ws.Range("F8").Value = 7.402862500103339
$ws.Range("G8").Value = 0.002812819207504053
$ws.Range("I8").Value = 5.536697672433931
$ws.Range("J8").Value = 0.1651213819755526
$ws.Range("L8").Value = 0.1599835084042525
$ws.Range("M8").Value = 12.447563431464

$ws.Range("C9").Value = 0.0629413649668038
$ws.Range("D9").Value = 0.1484523175404178
$ws.Range("E9").Value = 0.09729492870617662
$ws.Range("F9").Value = 6.9624230663307
$ws.Range("G9").Value = 0.002776609481830001
$ws.Range("I9").Value = 5.210013160579052
$ws.Range("J9").Value = 0.1792555236583695
$ws.Range("L9").Value = 0.1992755667670139
$ws.Range("M9").Value = 15.40976150285911

$ws.Range("C10").Value = 0.07331339805311643
$ws.Range("D10").Value = 0.161103748367708
$ws.Range("E10").Value = 0.1141967527032932
$ws.Range("F10").Value = 6.689803322038415
$ws.Range("G10").Value = 0.002751852736617974
$ws.Range("I10").Value = 5.00852890086415
$ws.Range("J10").Value = 0.1906892919340351
$ws.Range("L10").Value = 0.2285628323625986
$ws.Range("M10").Value = 17.59012777263155

$ws.Range("C11").Value = 0.07810358047638033
$ws.Range("D11").Value = 0.1671371359386455
$ws.Range("E11").Value = 0.1219076831428367
$ws.Range("F11").Value = 6.577337877566322
$ws.Range("G11").Value = 0.002740978076939132
$ws.Range("I11").Value = 4.925610423829681
$ws.Range("J11").Value = 0.1961424235683467
$ws.Range("L11").Value = 0.2419991677643623
$ws.Range("M11").Value = 18.58416102513786

$ws.Range("C12").Value = 0.07992879057482583
$ws.Range("D12").Value = 0.1694641621954815
$ws.Range("E12").Value = 0.1248315742270307
$ws.Range("F12").Value = 6.536454046529798
$ws.Range("G12").Value = 0.002736914721665212
$ws.Range("I12").Value = 4.895500809407508
$ws.Range("J12").Value = 0.1982457109846507
$ws.Range("L12").Value = 0.2471053718840324
$ws.Range("M12").Value = 18.96099667253446

$ws.Range("C13").Value = 0.07953518220909928
$ws.Range("D13").Value = 0.1689610755557283
$ws.Range("E13").Value = 0.1242016748873027
$ws.Range("F13").Value = 6.545182594569837
$ws.Range("G13").Value = 0.002737787424980957
$ws.Range("I13").Value = 4.901927541230393
$ws.Range("J13").Value = 0.1977909919202716
$ws.Range("L13").Value = 0.2460048210305388
$ws.Range("M13").Value = 18.87981808709748

$ws.Range("C14").Value = 0.07825351017990556
$ws.Range("D14").Value = 0.1673277201285543
$ws.Range("E14").Value = 0.1221481502576367
$ws.Range("F14").Value = 6.573939925657413
$ws.Range("G14").Value = 0.002740642693125144
$ws.Range("I14").Value = 4.923107241419203
$ws.Range("J14").Value = 0.1963146821038322
$ws.Range("L14").Value = 0.2424188832906395
$ws.Range("M14").Value = 18.61515454043945

$ws.Range("C15").Value = 0.0774699455408836
$ws.Range("D15").Value = 0.1663328237106043
$ws.Range("E15").Value = 0.1208908416286008
$ws.Range("F15").Value = 6.591777906826763
$ws.Range("G15").Value = 0.00274239871172373
$ws.Range("I15").Value = 4.936249407175623
$ws.Range("J15").Value = 0.1954154526864329
$ws.Range("L15").Value = 0.2402248145561146
$ws.Range("M15").Value = 18.45309799321001

$ws.Range("C16").Value = 0.07300187377170175
$ws.Range("D16").Value = 0.1607152450863509
$ws.Range("E16").Value = 0.1136933337128454
$ws.Range("F16").Value = 6.697389298214318
$ws.Range("G16").Value = 0.002752571125935191
$ws.Range("I16").Value = 5.014126458600231
$ws.Range("J16").Value = 0.1903381600355232
$ws.Range("L16").Value = 0.2276871729943082
$ws.Range("M16").Value = 17.52521782892046

$ws.Range("C17").Value = 0.07027991937536626
$ws.Range("D17").Value = 0.1573419018733944
$ws.Range("E17").Value = 0.1092840889070459
$ws.Range("F17").Value = 6.765167134117974
$ws.Range("G17").Value = 0.002758910036677581
$ws.Range("I17").Value = 5.064162688509327
$ws.Range("J17").Value = 0.187289350442768
$ws.Range("L17").Value = 0.2200260525067819
$ws.Range("M17").Value = 16.95662176432393

$ws.Range("C18").Value = 0.06872103137293095
$ws.Range("D18").Value = 0.1554276836887709
$ws.Range("E18").Value = 0.1067500505917707
$ws.Range("F18").Value = 6.80523564007521
$ws.Range("G18").Value = 0.002762592531388163
$ws.Range("I18").Value = 5.093762556633706
$ws.Range("J18").Value = 0.1855593275800231
$ws.Range("L18").Value = 0.2156301343311071
$ws.Range("M18").Value = 16.62977822053313

$ws.Range("C19").Value = 0.06819434313440809
$ws.Range("D19").Value = 0.1547839660564989
$ws.Range("E19").Value = 0.1058923980943192
$ws.Range("F19").Value = 6.818987118742001
$ws.Range("G19").Value = 0.002763845664268935
$ws.Range("I19").Value = 5.103924468676936
$ws.Range("J19").Value = 0.1849775588171667
$ws.Range("L19").Value = 0.2141435109746936
$ws.Range("M19").Value = 16.51914555310054

$ws.Range("C20").Value = 0.07056897518415894
$ws.Range("D20").Value = 0.1576982886391818
$ws.Range("E20").Value = 0.1097532439849687
$ws.Range("F20").Value = 6.757839519303673
$ws.Range("G20").Value = 0.002758231476483541
$ws.Range("I20").Value = 5.058751106506449
$ws.Range("J20").Value = 0.1876114471842811
$ws.Range("L20").Value = 0.2208404864629898
$ws.Range("M20").Value = 17.01712857753012

$ws.Range("C21").Value = 0.07862965472233441
$ws.Range("D21").Value = 0.1678063088641011
$ws.Range("E21").Value = 0.1227512074353001
$ws.Range("F21").Value = 6.565446587564082
$ws.Range("G21").Value = 0.002739802557483628
$ws.Range("I21").Value = 4.91685097145465
$ws.Range("J21").Value = 0.1967472531083558
$ws.Range("L21").Value = 0.2434716519536835
$ws.Range("M21").Value = 18.69288045546023

$ws.Range("C22").Value = 0.08396394512249117
$ws.Range("D22").Value = 0.1746602738161585
$ws.Range("E22").Value = 0.1312694126686793
$ws.Range("F22").Value = 6.449662824021971
$ws.Range("G22").Value = 0.002728076171437015
$ws.Range("I22").Value = 4.831645663340552
$ws.Range("J22").Value = 0.2029423787512883
$ws.Range("L22").Value = 0.2583692132072031
$ws.Range("M22").Value = 19.79056508420888

$ws.Range("C23").Value = 0.08111055801740008
$ws.Range("D23").Value = 0.1709787176865802
$ws.Range("E23").Value = 0.1267206974378823
$ws.Range("F23").Value = 6.510532925389384
$ws.Range("G23").Value = 0.002734306028558171
$ws.Range("I23").Value = 4.876420491868444
$ws.Range("J23").Value = 0.1996146659907225
$ws.Range("L23").Value = 0.2504076850617878
$ws.Range("M23").Value = 19.20444626355106

$ws.Range("C24").Value = 0.07043827443524719
$ws.Range("D24").Value = 0.1575370879828313
$ws.Range("E24").Value = 0.1095411362406935
$ws.Range("F24").Value = 6.761148904348772
$ws.Range("G24").Value = 0.002758538134469191
$ws.Range("I24").Value = 5.061195089463212
$ws.Range("J24").Value = 0.1874657564541025
$ws.Range("L24").Value = 0.2204722543407485
$ws.Range("M24").Value = 16.98977329961127

$ws.Range("C25").Value = 0.0591795324558575
$ws.Range("D25").Value = 0.1440251814783835
$ws.Range("E25").Value = 0.09108573358493999
$ws.Range("F25").Value = 7.07278214119259
$ws.Range("G25").Value = 0.002786076605364353
$ws.Range("I25").Value = 5.291746182391094
$ws.Range("J25").Value = 0.1752551375975315
$ws.Range("L25").Value = 0.2460048210305388
$ws.Range("M25").Value = 18.87981808709748
